# Update "想去人数" (people interested) counts for two bilibili events
# that appear on both the "展览" sheet and the consolidated "全部类型" sheet.
#   id=91123 event: 420 -> 421
#   id=92565 event: 2667 -> 2691

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 421
$ws1.Range("F3").Value = 2691

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 421
$ws4.Range("F7").Value = 2691
